# Update the "Country WHS" article-number codes on the Template sheet.
# The target cells hold short, zero-padded numeric codes (e.g. "05") that
# must remain text (so leading zeros survive) while keeping the cells'
# existing numeric-looking style. We briefly switch the cell to a Text
# number format to enter the value as a string, then restore the original
# "0" number format - this keeps reusing the workbook's existing style
# (no new styles / quote-prefixed cells get introduced).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Template")

function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.NumberFormat = "0"
}

Set-TextValue $ws.Range("C2") "30"
Set-TextValue $ws.Range("C3") "32"
Set-TextValue $ws.Range("C4") "65"
Set-TextValue $ws.Range("C5") "05"
